$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.252.61"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "1.557.09"
$ws.Range("E3").Value = "  -4.04%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.17"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.475"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0602"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.240"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.69"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "1.776.03"
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").Value = "1.568.29"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.97"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.500"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").Value = "25.196.38"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.78"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "184.52"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.08"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.21"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.17"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.65"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.65%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.74"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.39"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.72%  "
$ws.Range("E30").Value = "  -6.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0461"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.00"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.96"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.44"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("E35").Value = "  -4.01%  "
$ws.Range("D36").Value = "1.083.79"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.54%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0149"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.815"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.490"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.758"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.33"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.04"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "1.690.51"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.13"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("E48").Value = "  -4.93%  "
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("E51").Value = "  -0.29%  "
